$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Variables")
$ws2 = $wb.Worksheets.Item("Categories")

# --- Sheet1 ("Variables") ---
# Rename B36: inccanc_fup5 -> inccanc (label/type for that row stay as-is)
$ws1.Range("B36").Value = "inccanc"

# The old row 41 ("inccanc" / incident first occuring cancer / integer) is
# removed; rows 42..79's name/label/valueType (columns B:D) shift up into
# rows 41..78, while column A (the numeric index) stays fixed per row.
for ($r = 41; $r -le 78; $r++) {
    $ws1.Cells.Item($r, 2).Value = $ws1.Cells.Item($r + 1, 2).Value2
    $ws1.Cells.Item($r, 3).Value = $ws1.Cells.Item($r + 1, 3).Value2
    $ws1.Cells.Item($r, 4).Value = $ws1.Cells.Item($r + 1, 4).Value2
}

# Fix the (now shifted) label text at row 41
$ws1.Range("C41").Value = "censored date for cancer at FUP5"

# Drop the now-duplicated trailing row entirely, shrinking the dimension
$ws1.Rows("79").Delete()

# --- Sheet2 ("Categories") ---
# Remove the 3 lookup rows that belonged to the retired "inccanc_fup5"
# variable; everything below shifts up to close the gap.
$ws2.Rows("67:69").Delete()
